$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-28 Thursday" "2024-11-29 Friday"

Replace-Text "63×57=" "42×30="
Replace-Text "64×81=" "19×35="
Replace-Text "21×57=" "35×58="
Replace-Text "60×59=" "39×73="
Replace-Text "86×96=" "58×37="
Replace-Text "53×37=" "90×38="
Replace-Text "97×67=" "38×80="
Replace-Text "67×67=" "95×56="
Replace-Text "71×60=" "74×22="
Replace-Text "11×28=" "78×94="
Replace-Text "96×18=" "90×33="
Replace-Text "48×34=" "94×41="
Replace-Text "91×25=" "72×81="
Replace-Text "47×68=" "81×14="
Replace-Text "45×48=" "54×91="
Replace-Text "90×92=" "31×69="
Replace-Text "52×25=" "57×22="
Replace-Text "31×99=" "90×97="
Replace-Text "68×84=" "55×67="
Replace-Text "21×58=" "85×29="
Replace-Text "37×47=" "70×60="
Replace-Text "66×11=" "52×76="
Replace-Text "72×42=" "87×25="
Replace-Text "21×20=" "11×16="
Replace-Text "37×90=" "19×82="
